$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4330042.5
$ws.Range("I74").Value = 6492064
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 6492064
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -6491128
$ws.Range("N74").Value = -7872

$ws.Range("H77").Value = 4330042.5
$ws.Range("I77").Value = 6492064
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 32460320
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -32455640
$ws.Range("N77").Value = -39360

$ws.Range("H86").Value = 2414.261
$ws.Range("I86").Value = 1295.2
$ws.Range("K86").Value = 1295.2
$ws.Range("M86").Value = -172.2

$ws.Range("H88").Value = 662930.1
$ws.Range("I88").Value = 2300.4285
$ws.Range("J88").Value = 934954.1
$ws.Range("K88").Value = 2300.4285
$ws.Range("L88").Value = 934954.1
$ws.Range("M88").Value = -1894.4285
$ws.Range("N88").Value = -935766.1

$ws.Range("H89").Value = 2414.261
$ws.Range("I89").Value = 1295.2
$ws.Range("K89").Value = 6476
$ws.Range("M89").Value = -860

$ws.Range("H91").Value = 662930.1
$ws.Range("I91").Value = 2300.4285
$ws.Range("J91").Value = 934954.1
$ws.Range("K91").Value = 2300.4285
$ws.Range("L91").Value = 934954.1
$ws.Range("M91").Value = -896.4285
$ws.Range("N91").Value = -937762.1

$ws.Range("H125").Value = 2868.8948
$ws.Range("I125").Value = 1059.3529
$ws.Range("J125").Value = 18250
$ws.Range("K125").Value = 9534.176100000001
$ws.Range("L125").Value = 164250
$ws.Range("M125").Value = -7074.176100000001
$ws.Range("N125").Value = -169170

$ws.Range("H137").Value = 835.2692
$ws.Range("I137").Value = 772.0952
$ws.Range("J137").Value = 1100.6
$ws.Range("K137").Value = 2316.2856
$ws.Range("L137").Value = 3301.8
$ws.Range("M137").Value = 233.7143999999998
$ws.Range("N137").Value = -8401.799999999999

$ws.Range("H138").Value = 2987.01
$ws.Range("I138").Value = 727.75
$ws.Range("J138").Value = 3700.4604
$ws.Range("K138").Value = 2183.25
$ws.Range("L138").Value = 11101.3812
$ws.Range("M138").Value = 2956.75
$ws.Range("N138").Value = -21381.3812

$ws.Range("H141").Value = 3136.923
$ws.Range("I141").Value = 2343.6365
$ws.Range("J141").Value = 7500
$ws.Range("K141").Value = 7030.9095
$ws.Range("L141").Value = 22500
$ws.Range("M141").Value = -1850.9095
$ws.Range("N141").Value = -32860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1000
$ws.Range("J26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("N26").Value = -1660

$ws.Range("H31").Value = 30000
$ws.Range("I31").Value = 8000
$ws.Range("J31").Value = 52000
$ws.Range("K31").Value = 8000
$ws.Range("L31").Value = 52000
$ws.Range("M31").Value = -7706
$ws.Range("N31").Value = -52588

$ws.Range("H32").Value = 5273.623
$ws.Range("I32").Value = 4396.98
$ws.Range("K32").Value = 4396.98
$ws.Range("M32").Value = -4109.98

$ws.Range("H33").Value = 1026
$ws.Range("I33").Value = 1026
$ws.Range("K33").Value = 1026
$ws.Range("M33").Value = -697

$ws.Range("H88").Value = 2757.8
$ws.Range("J88").Value = 2837
$ws.Range("L88").Value = 2837
$ws.Range("N88").Value = -3649

$ws.Range("H91").Value = 2757.8
$ws.Range("J91").Value = 2837
$ws.Range("L91").Value = 2837
$ws.Range("N91").Value = -5645

$ws.Range("H102").Value = 4223.077
$ws.Range("I102").Value = 4000
$ws.Range("J102").Value = 4725
$ws.Range("K102").Value = 4000
$ws.Range("L102").Value = 4725
$ws.Range("M102").Value = -2378
$ws.Range("N102").Value = -7969

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1853.25
$ws.Range("I86").Value = 1803.8235
$ws.Range("J86").Value = 2133.3333
$ws.Range("K86").Value = 1803.8235
$ws.Range("L86").Value = 2133.3333
$ws.Range("M86").Value = -680.8235
$ws.Range("N86").Value = -4379.3333

$ws.Range("H89").Value = 1853.25
$ws.Range("I89").Value = 1803.8235
$ws.Range("J89").Value = 2133.3333
$ws.Range("K89").Value = 9019.1175
$ws.Range("L89").Value = 10666.6665
$ws.Range("M89").Value = -3403.1175
$ws.Range("N89").Value = -21898.6665

$ws.Range("H103").Value = 18328.5
$ws.Range("J103").Value = 18328.5
$ws.Range("L103").Value = 18328.5
$ws.Range("N103").Value = -20672.5

$ws.Range("H105").Value = 4929.07
$ws.Range("I105").Value = 4814.0938
$ws.Range("J105").Value = 5263.5454
$ws.Range("K105").Value = 4814.0938
$ws.Range("L105").Value = 5263.5454
$ws.Range("M105").Value = -3067.0938
$ws.Range("N105").Value = -8757.545399999999

$ws.Range("H134").Value = 36268.1
$ws.Range("I134").Value = 1620.625
$ws.Range("J134").Value = 202576
$ws.Range("K134").Value = 4861.875
$ws.Range("L134").Value = 607728
$ws.Range("M134").Value = -2326.875
$ws.Range("N134").Value = -612798

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 6000000
$ws.Range("I32").Value = 6000000
$ws.Range("K32").Value = 6000000
$ws.Range("M32").Value = -5999684

$ws.Range("H62").Value = 8594.736999999999
$ws.Range("I62").Value = 8883.333000000001
$ws.Range("J62").Value = 8100
$ws.Range("K62").Value = 8883.333000000001
$ws.Range("L62").Value = 8100
$ws.Range("M62").Value = -8259.333000000001
$ws.Range("N62").Value = -9348

$ws.Range("H65").Value = 8594.736999999999
$ws.Range("I65").Value = 8883.333000000001
$ws.Range("J65").Value = 8100
$ws.Range("K65").Value = 44416.665
$ws.Range("L65").Value = 40500
$ws.Range("M65").Value = -41296.665
$ws.Range("N65").Value = -46740

$ws.Range("H105").Value = 800
$ws.Range("I105").Value = 800
$ws.Range("K105").Value = 800
$ws.Range("M105").Value = 947

$ws.Range("H132").Value = 2952
$ws.Range("I132").Value = 2109
$ws.Range("J132").Value = 3888.6667
$ws.Range("K132").Value = 6327
$ws.Range("L132").Value = 11666.0001
$ws.Range("M132").Value = -3797
$ws.Range("N132").Value = -16726.0001

$ws.Range("H134").Value = 2556.7778
$ws.Range("I134").Value = 1843.0416
$ws.Range("J134").Value = 8266.666999999999
$ws.Range("K134").Value = 5529.1248
$ws.Range("L134").Value = 24800.001
$ws.Range("M134").Value = -2994.1248
$ws.Range("N134").Value = -29870.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4425.6665
$ws.Range("I70").Value = 4079.9443
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 4079.9443
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -3809.9443
$ws.Range("N70").Value = -7040

$ws.Range("H73").Value = 4425.6665
$ws.Range("I73").Value = 4079.9443
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 4079.9443
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = -3143.9443
$ws.Range("N73").Value = -8372

$ws.Range("H132").Value = 2390.95
$ws.Range("I132").Value = 1838.6923
$ws.Range("J132").Value = 3416.5715
$ws.Range("K132").Value = 5516.0769
$ws.Range("L132").Value = 10249.7145
$ws.Range("M132").Value = -2986.0769
$ws.Range("N132").Value = -15309.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1700.6
$ws.Range("I100").Value = 1434.3334
$ws.Range("K100").Value = 1434.3334
$ws.Range("M100").Value = -893.3334

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 2604.796
$ws.Range("I132").Value = 2459.5386
$ws.Range("J132").Value = 3171.3
$ws.Range("K132").Value = 7378.6158
$ws.Range("L132").Value = 9513.900000000001
$ws.Range("M132").Value = -4848.6158
$ws.Range("N132").Value = -14573.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2215.1538
$ws.Range("I132").Value = 1713.4286
$ws.Range("J132").Value = 2800.5
$ws.Range("K132").Value = 5140.2858
$ws.Range("L132").Value = 8401.5
$ws.Range("M132").Value = -2610.2858
$ws.Range("N132").Value = -13461.5

$ws.Range("H136").Value = 2312.2285
$ws.Range("I136").Value = 2242.4167
$ws.Range("J136").Value = 2464.5454
$ws.Range("K136").Value = 6727.250100000001
$ws.Range("L136").Value = 7393.6362
$ws.Range("M136").Value = -4177.250100000001
$ws.Range("N136").Value = -12493.6362
